# Apply the changes described by the diff:
#  - Metadata!B8 (Date value): 2025-07-18T06:40:38+00:00 -> 2025-07-21T11:52:46+00:00
#  - Include #0!B4 (System URI): ansforge.github.io -> interop.esante.gouv.fr (TRE-G15-ProfessionSante)
#  - Include #1!B4 (System URI): ansforge.github.io -> interop.esante.gouv.fr (TRE-R95-UsagerTitre)
#  - Include #2!B4 (System URI): ansforge.github.io -> interop.esante.gouv.fr (TRE-R94-ProfessionSocial)
#  - Include #3!B4 (System URI): ansforge.github.io -> interop.esante.gouv.fr (TRE-R291-AutreProfession)

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B8").Value = "2025-07-21T11:52:46+00:00"

$wsInclude0 = $wb.Worksheets.Item("Include #0")
$wsInclude0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-G15-ProfessionSante"

$wsInclude1 = $wb.Worksheets.Item("Include #1")
$wsInclude1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R95-UsagerTitre"

$wsInclude2 = $wb.Worksheets.Item("Include #2")
$wsInclude2.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R94-ProfessionSocial"

$wsInclude3 = $wb.Worksheets.Item("Include #3")
$wsInclude3.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R291-AutreProfession"
